# Append a new data row (row 17) to the case-tracking sheet, mirroring the
# existing rows (A..I populated, J..M left blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Range("A$row").Value = "YXYZH"
# Leading apostrophe forces this to be stored as literal text rather than
# being auto-parsed into a date serial number (matches the existing rows,
# which store dates as plain text strings like "11/29/2024").
$ws.Range("B$row").Value = "'11/29/2024"
$ws.Range("C$row").Value = "Ashar Nadeem"
$ws.Range("D$row").Value = "0322-7287568"
$ws.Range("E$row").Value = "Lahore"
$ws.Range("F$row").Value = "0322-7287568"
$ws.Range("G$row").Value = "xxdxsdxdsx"
$ws.Range("H$row").Value = "Iqbal town"
$ws.Range("I$row").Value = "knj"

# Reset to the default "Normal" style so the new cells don't pick up any
# incidental formatting (e.g. the quote-prefix flag from the literal-text
# trick above, or an explicitly-stamped column style) and instead match the
# unstyled cells used by the other data rows.
$ws.Range("A$($row):I$row").Style = "Normal"
